$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categóricos")

# Rename the three "codigo_*" header cells to their shorter "cve*" forms.
$ws.Range("A2").Value = "cvevisita"
$ws.Range("A10").Value = "cvecondicion"
$ws.Range("A18").Value = "cvepais"

$ws.Range("A18").Select()
